$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("W3").Value = "[327, 615, 626, 786, 842, 863, 950, 973, 1040, 1068, 1130]"
$ws.Range("Y3").Value = "{'initial_window': 254, 'window_size': 241, 'num_trees': 48, 'max_depth': 19}"
$ws.Range("Z3").Value = 466

# Row 4 updates
$ws.Range("W4").Value = "[563]"
$ws.Range("Y4").Value = "{'initial_window': 141, 'window_size': 236, 'num_trees': 43, 'max_depth': 18}"
$ws.Range("Z4").Value = 215

# Row 6 updates
$ws.Range("W6").Value = "[579, 990]"
$ws.Range("X6").Value = 1
$ws.Range("Y6").Value = "{'initial_window': 219, 'window_size': 203, 'num_trees': 31, 'max_depth': 14}"
$ws.Range("Z6").Value = 8
